# Fixed naive component forecaster bug - Presentation state 11.02.
# For each data row (2-16), shift existing values one column to the right
# (processed right-to-left to avoid clobbering unread values), then write
# the new forecast-error value into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = -0.5607856406117955
$ws.Range("J2").Value = -0.06012564061179543
$ws.Range("I2").Value = -0.2731016406117954
$ws.Range("H2").Value = -0.3657496406117954
$ws.Range("G2").Value = -0.3937486406117954
$ws.Range("F2").Value = -0.1063816406117954
$ws.Range("E2").Value = -0.6170576406117954
$ws.Range("D2").Value = -0.3099216406117954
$ws.Range("C2").Value = -0.5449296406117954
$ws.Range("B2").Value = 0.3847923593882046

$ws.Range("K3").Value = 0.3115658964218456
$ws.Range("J3").Value = 0.02180889642184558
$ws.Range("I3").Value = 0.5224688964218456
$ws.Range("H3").Value = 0.3094928964218456
$ws.Range("G3").Value = 0.2168448964218456
$ws.Range("F3").Value = 0.1888458964218456
$ws.Range("E3").Value = 0.4762128964218456
$ws.Range("D3").Value = -0.03446310357815441
$ws.Range("C3").Value = 0.2726728964218456
$ws.Range("B3").Value = 0.03766489642184559

$ws.Range("K4").Value = -0.04200114438751221
$ws.Range("J4").Value = 0.2158908556124878
$ws.Range("I4").Value = -0.07386614438751221
$ws.Range("H4").Value = 0.4267938556124878
$ws.Range("G4").Value = 0.2138178556124878
$ws.Range("F4").Value = 0.1211698556124878
$ws.Range("E4").Value = 0.09317085561248779
$ws.Range("D4").Value = 0.3805378556124878
$ws.Range("C4").Value = -0.1301381443875122
$ws.Range("B4").Value = 0.1769978556124878

$ws.Range("K5").Value = 0.6007988926112107
$ws.Range("J5").Value = 0.2184178926112106
$ws.Range("I5").Value = 0.4763098926112106
$ws.Range("H5").Value = 0.1865528926112106
$ws.Range("G5").Value = 0.6872128926112107
$ws.Range("F5").Value = 0.4742368926112106
$ws.Range("E5").Value = 0.3815888926112106
$ws.Range("D5").Value = 0.3535898926112106
$ws.Range("C5").Value = 0.6409568926112106
$ws.Range("B5").Value = 0.1302808926112106

$ws.Range("K6").Value = -0.6459331975472806
$ws.Range("J6").Value = -0.2346561035472806
$ws.Range("I6").Value = -0.6170371035472806
$ws.Range("H6").Value = -0.3591451035472806
$ws.Range("G6").Value = -0.6489021035472806
$ws.Range("F6").Value = -0.1482421035472806
$ws.Range("E6").Value = -0.3612181035472806
$ws.Range("D6").Value = -0.4538661035472806
$ws.Range("C6").Value = -0.4818651035472806
$ws.Range("B6").Value = -0.1944981035472806

$ws.Range("K7").Value = -0.2588977436446591
$ws.Range("J7").Value = -0.5458538376446591
$ws.Range("I7").Value = -0.1345767436446591
$ws.Range("H7").Value = -0.5169577436446591
$ws.Range("G7").Value = -0.2590657436446591
$ws.Range("F7").Value = -0.5488227436446591
$ws.Range("E7").Value = -0.04816274364465911
$ws.Range("D7").Value = -0.2611387436446591
$ws.Range("C7").Value = -0.3537867436446591
$ws.Range("B7").Value = -0.3817857436446591

$ws.Range("J8").Value = 0.09488876243503713
$ws.Range("I8").Value = -0.1920673315649629
$ws.Range("H8").Value = 0.2192097624350371
$ws.Range("G8").Value = -0.1631712375649629
$ws.Range("F8").Value = 0.09472076243503715
$ws.Range("E8").Value = -0.1950362375649629
$ws.Range("D8").Value = 0.3056237624350371
$ws.Range("C8").Value = 0.09264776243503714
$ws.Range("B8").Value = -0.0000002375649628613696

$ws.Range("I9").Value = -0.0532886279974082
$ws.Range("H9").Value = -0.3402447219974082
$ws.Range("G9").Value = 0.0710323720025918
$ws.Range("F9").Value = -0.3113486279974082
$ws.Range("E9").Value = -0.05345662799740819
$ws.Range("D9").Value = -0.3432136279974082
$ws.Range("C9").Value = 0.1574463720025918
$ws.Range("B9").Value = -0.0555296279974082

$ws.Range("H10").Value = -0.2107346170015632
$ws.Range("G10").Value = -0.4976907110015631
$ws.Range("F10").Value = -0.08641361700156319
$ws.Range("E10").Value = -0.4687946170015632
$ws.Range("D10").Value = -0.2109026170015632
$ws.Range("C10").Value = -0.5006596170015631
$ws.Range("B10").Value = 0.0000003829984367986761

$ws.Range("G11").Value = 0.2899248395245076
$ws.Range("F11").Value = 0.002968745524507627
$ws.Range("E11").Value = 0.4142458395245076
$ws.Range("D11").Value = 0.03186483952450761
$ws.Range("C11").Value = 0.2897568395245076
$ws.Range("B11").Value = -0.0000001604754923945073

$ws.Range("F12").Value = 0.009566958989038449
$ws.Range("E12").Value = -0.2773891350109615
$ws.Range("D12").Value = 0.1338879589890384
$ws.Range("C12").Value = -0.2484930410109615
$ws.Range("B12").Value = 0.009398958989038461

$ws.Range("E13").Value = 0.1815418168340501
$ws.Range("D13").Value = -0.1054142771659499
$ws.Range("C13").Value = 0.3058628168340501
$ws.Range("B13").Value = -0.07651818316594991

$ws.Range("D14").Value = -0.1243207229122814
$ws.Range("C14").Value = -0.4112768169122814
$ws.Range("B14").Value = 0.0000002770877186031306

$ws.Range("C15").Value = 0.4880092297750048
$ws.Range("B15").Value = 0.2010531357750048

$ws.Range("B16").Value = -0.2003621554241067
